$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 9, shifting existing rows 9-14 down to 11-16
$ws.Rows.Item(9).Resize(2).EntireRow.Insert()

# Column A holds date-like text (e.g. "2023-11-07") that Excel would
# otherwise auto-recognize and silently convert into a date serial number
# with a date number-format style. Force the cells to Text first so the
# literal string is preserved (this was the "styles" bug being fixed),
# then the stray auto-applied style is cleared below once the values are in.
$ws.Range("A9:A10").NumberFormat = "@"
$ws.Range("A17:A18").NumberFormat = "@"

# Fill the two newly inserted rows (9 and 10) with new data
$ws.Range("A9").Value = "2023-11-07"
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = "expenses"
$ws.Range("D9").Value = 93.04000000000001
$ws.Range("E9").Value = 99.01000000000001
$ws.Range("F9").Value = 12.69

$ws.Range("A10").Value = "2023-11-10"
$ws.Range("B10").Value = 250
$ws.Range("C10").Value = "expenses"
$ws.Range("D10").Value = 91.93000000000001
$ws.Range("E10").Value = 98.41
$ws.Range("F10").Value = 12.6

# Append two new rows (17 and 18) with new data
$ws.Range("A17").Value = "2023-11-01"
$ws.Range("B17").Value = 25
$ws.Range("C17").Value = "income"
$ws.Range("D17").Value = 92.02
$ws.Range("E17").Value = 97.93000000000001
$ws.Range("F17").Value = 12.53

$ws.Range("A18").Value = "2023-11-11"
$ws.Range("B18").Value = 1000
$ws.Range("C18").Value = "income"
$ws.Range("D18").Value = 92.05
$ws.Range("E18").Value = 98.31999999999999
$ws.Range("F18").Value = 12.6

# Now that the text values are locked in, clear the cell formatting so the
# cells fall back to the workbook's default (unstyled) cell, matching every
# other data cell in the sheet (no stray number-format style left behind).
$ws.Range("A9:A10").ClearFormats()
$ws.Range("A17:A18").ClearFormats()
